$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 / Tableau1: add a calculated "Formula" column ---
$lo = $ws1.ListObjects.Item("Tableau1")
$col = $lo.ListColumns.Add()

# Header (also renames the table column + registers the shared string)
$ws1.Range("E1").Value = "Formula"

# Calculated column formula (per data row, table-style structured reference)
$ws1.Range("E2").Formula = "=Tableau1[[#This Row],[Integers]]*Tableau1[[#This Row],[Floats]]"
$ws1.Range("E3").Formula = "=Tableau1[[#This Row],[Integers]]*Tableau1[[#This Row],[Floats]]"
$ws1.Range("E4").Formula = "=Tableau1[[#This Row],[Integers]]*Tableau1[[#This Row],[Floats]]"

# Keep the new data cells on the default "Normal" style (matches the
# unformatted calculated column produced by the original edit).
$ws1.Range("E2:E4").Style = "Normal"

# --- Sheet2: move the active selection to D5 without changing the ---
# --- workbook's active sheet (Sheet1 stays the tab in focus)       ---
[void]$ws2.Range("D5").Select()
[void]$ws1.Activate()
